$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-mark the Price column (D) as Text so numeric-looking strings
# (e.g. "603.84") are stored as literal text, not coerced to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.301.61"
$ws.Range("E2").Value = "  +0.63%  "

# Row 3
$ws.Range("D3").Value = "3.510.54"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "603.84"
$ws.Range("E5").Value = "  +4.38%  "

# Row 6
$ws.Range("D6").Value = "169.99"
$ws.Range("E6").Value = "  -3.18%  "

# Row 7
$ws.Range("D7").Value = "0.614"
$ws.Range("E7").Value = "  -1.16%  "

# Row 8
$ws.Range("D8").Value = "3.505.23"
$ws.Range("E8").Value = "  +0.22%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "0.193"
$ws.Range("E10").Value = "  +2.02%  "

# Row 11
$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  +1.14%  "

# Row 12
$ws.Range("E12").Value = "  -3.76%  "

# Row 13
$ws.Range("D13").Value = "47.05"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
$ws.Range("D14").Value = "0.0000277"
$ws.Range("E14").Value = "  +0.00%  "

# Row 15
$ws.Range("D15").Value = "4.075.37"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16
$ws.Range("D16").Value = "8.33"
$ws.Range("E16").Value = "  -6.49%  "

# Row 17
$ws.Range("D17").Value = "612.23"
$ws.Range("E17").Value = "  -9.37%  "

# Row 18
$ws.Range("D18").Value = "3.512.04"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19
$ws.Range("D19").Value = "69.338.12"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
$ws.Range("E20").Value = "  -1.53%  "

# Row 21
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -1.89%  "

# Row 22
$ws.Range("D22").Value = "11.20"
$ws.Range("E22").Value = "  +0.40%  "

# Row 23
$ws.Range("D23").Value = "0.879"
$ws.Range("E23").Value = "  -2.89%  "

# Row 24
$ws.Range("D24").Value = "15.75"
$ws.Range("E24").Value = "  -3.47%  "

# Row 25
$ws.Range("D25").Value = "96.02"
$ws.Range("E25").Value = "  -2.27%  "

# Row 26
$ws.Range("D26").Value = "3.85"
$ws.Range("E26").Value = "  -0.35%  "

# Row 27
$ws.Range("E27").Value = "  -0.04%  "

# Row 28
$ws.Range("D28").Value = "2.61"
$ws.Range("E28").Value = "  -1.69%  "

# Row 29
$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  -2.00%  "

# Row 30
$ws.Range("D30").Value = "32.97"
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("D31").Value = "3.10"
$ws.Range("E31").Value = "  -3.66%  "

# Row 32
$ws.Range("D32").Value = "8.43"
$ws.Range("E32").Value = "  -3.76%  "

# Row 33
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  -1.23%  "

# Row 34
$ws.Range("D34").Value = "6.89"
$ws.Range("E34").Value = "  -6.54%  "

# Row 35
$ws.Range("D35").Value = "567.51"
$ws.Range("E35").Value = "  -0.22%  "

# Row 36
$ws.Range("D36").Value = "10.75"
$ws.Range("E36").Value = "  -1.87%  "

# Row 37
$ws.Range("D37").Value = "3.49"
$ws.Range("E37").Value = "  -3.87%  "

# Row 38
$ws.Range("D38").Value = "0.101"
$ws.Range("E38").Value = "  -3.98%  "

# Row 39
$ws.Range("D39").Value = "56.85"
$ws.Range("E39").Value = "  +0.17%  "

# Row 40
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0447"
$ws.Range("E41").Value = "  +1.36%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.139"
$ws.Range("E42").Value = "  +1.54%  "

# Row 43
$ws.Range("D43").Value = "3.370.56"
$ws.Range("E43").Value = "  -1.47%  "

# Row 44
$ws.Range("D44").Value = "0.325"
$ws.Range("E44").Value = "  -3.69%  "

# Row 45
$ws.Range("D45").Value = "32.93"
$ws.Range("E45").Value = "  -1.42%  "

# Row 46
$ws.Range("D46").Value = "0.0₃0699"
$ws.Range("E46").Value = "  -0.73%  "

# Row 47
$ws.Range("D47").Value = "2.89"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("E49").Value = "  -3.12%  "

# Row 50
$ws.Range("D50").Value = "134.47"
$ws.Range("E50").Value = "  +0.58%  "

# Row 51
$ws.Range("D51").Value = "5.69"
$ws.Range("E51").Value = "  +9.59%  "

# Remove the temporary Text format again so the cells end up with
# no explicit style, matching the original (style-less) cells.
$dRange.ClearFormats()